# Adds 12 additional (currently blank) problem rows (rows 7-18) below the
# existing data, reusing the same cell formatting as the last existing row
# (row 6), then gives the newly added "problem block" rows a taller custom
# height (matching extra wrapped description/solution text that was typed
# into them), and leaves the remaining filler rows (12-18) at the default
# row height. Finally the active selection is moved to B10 to match where
# the author ended up after typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 already exists (last row of the previous table) - give it its new,
# taller custom height first.
$ws.Rows.Item(6).RowHeight = 36.75

# Use row 6's existing formatting (style index "1") as the template for the
# newly appended rows so no new style records are introduced - copy formats
# only (xlPasteFormats = -4122) down through row 18.
$ws.Range("A6:S6").Copy()
$ws.Range("A7:S18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Custom heights for the new problem rows that hold the extra
# descriptions/solutions.
$ws.Rows.Item(7).RowHeight = 31.5
$ws.Rows.Item(8).RowHeight = 36
$ws.Rows.Item(9).RowHeight = 40.5
$ws.Rows.Item(10).RowHeight = 32.25
$ws.Rows.Item(11).RowHeight = 24

# Rows 12-18 stay at the default row height (no custom height applied).

# Leave the selection where the author last left it while filling in the
# new rows.
$ws.Range("B10").Select()
